# Refresh the "realeffort" column with new values and re-sort the ranking
# table (rows 2-13) in descending order of realeffort, recomputing the
# re_rank column (H) to match the new order. Each worker keeps her own
# index/prolificid/name/gender/race; only realeffort + re_rank + row order
# change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of prolificid -> updated realeffort value.
$newEffort = @{
    "60bfcf5805c5ae12a546f9f3" = 7.186737775649103
    "5c5882fc5bfe7600011197cb" = 6.098775820819761
    "60b45e9961dd412bfb6780f8" = 6.491841461046875
    "60bd88b8fc436774352f53b9" = 5.03843188797754
    "5e96194b0a9fe909389e9f7b" = 5.020761003118488
    "60c0e5899d387663c07eb3a4" = 4.239120605821088
    "6077db0613ce87b4a62a78f9" = 1.248319253184411
    "608b14a312c099ac00b721b6" = 1.073514296423548
    "60cb36ee9f58331a33cf5506" = 0.1805031549777598
    "60d5775a99b502eec8cf56b4" = 0.3953773248513697
    "6036f9b3b1842f8b659b18c7" = 0.4217658868825644
    "5c0e89c6c323400001e6c4a5" = 0.4004490700212808
}

$firstRow = 2
$lastRow = 13

# Read the existing rows (prolificid, name, gender, race) and apply the
# new realeffort value per worker.
$rows = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $prolificid = $ws.Cells.Item($r, 3).Value2
    $row = [PSCustomObject]@{
        Index       = $ws.Cells.Item($r, 2).Value2
        Prolificid  = $prolificid
        Name        = $ws.Cells.Item($r, 4).Value2
        Gender      = $ws.Cells.Item($r, 5).Value2
        Realeffort  = $newEffort[$prolificid]
        Race        = $ws.Cells.Item($r, 7).Value2
    }
    $rows += $row
}

# Sort by the refreshed realeffort value, descending.
$sorted = $rows | Sort-Object -Property Realeffort -Descending

# Write the re-sorted rows back, recomputing re_rank (column H).
$r = $firstRow
$rank = 1
foreach ($row in $sorted) {
    $ws.Cells.Item($r, 2).Value = $row.Index
    $ws.Cells.Item($r, 3).Value = $row.Prolificid
    $ws.Cells.Item($r, 4).Value = $row.Name
    $ws.Cells.Item($r, 5).Value = $row.Gender
    $ws.Cells.Item($r, 6).Value = $row.Realeffort
    $ws.Cells.Item($r, 7).Value = $row.Race
    $ws.Cells.Item($r, 8).Value = $rank
    $r++
    $rank++
}
